# Generate Report for handoff
# Adds two newly-tracked source files (32217ea7... and bea9fae3...) to the
# localization status report: one row per file in the Overview sheet, plus
# one row per file in each language sheet (zh-cn, de-de). Also flips the
# status of the two files that already had a handoff generated previously
# (c2b13a0c... and dbf6e8a4...) from "Ready for handoff" to "In Translation".

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/e2b02ab458a2a4cedd27d5d51f652ddd8809a1a6/e2e/"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/e2b02ab458a2a4cedd27d5d51f652ddd8809a1a6/.localization-config"
$xlfZh    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/00f28509a6c1732639fd524a5e39dbe9ddbb2a09/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$xlfDe    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35d88ba1f854fbdbbf271e0659bc387e9b4bd88e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

$file1md  = "32217ea7-4110-42d9-960c-31dcc762cd2a.md"
$file2md  = "bea9fae3-5e2b-487a-8252-b0cdcb1f2603.md"

$file1xlfZh = "32217ea7-4110-42d9-960c-31dcc762cd2a.36a380d573f39b0e6c80450d3743ffe756a1aa54.zh-cn.xlf"
$file2xlfZh = "bea9fae3-5e2b-487a-8252-b0cdcb1f2603.3c91257f2c42de64eac33bbaf22ad3205513821c.zh-cn.xlf"
$file1xlfDe = "32217ea7-4110-42d9-960c-31dcc762cd2a.36a380d573f39b0e6c80450d3743ffe756a1aa54.de-de.xlf"
$file2xlfDe = "bea9fae3-5e2b-487a-8252-b0cdcb1f2603.3c91257f2c42de64eac33bbaf22ad3205513821c.de-de.xlf"

$handoffZh = "2016-02-16 14:28:37"
$handoffDe = "2016-02-16 14:28:51"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Rows for the two already-handed-off files move from "Ready for handoff"
# to "In Translation".
$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"
$ov.Range("B3").Value = "In Translation"
$ov.Range("C3").Value = "In Translation"

# Row 4 used to be the ".localization-config" row; it now becomes the first
# newly tracked file, and the config row is pushed down to row 6.
$ov.Hyperlinks.Add($ov.Range("A4"), ($mdBase + $file1md), "", "", $file1md)
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"

$ov.Hyperlinks.Add($ov.Range("A5"), ($mdBase + $file2md), "", "", $file2md)
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"

$ov.Hyperlinks.Add($ov.Range("A6"), $cfgUrl, "", "", ".localization-config")
$ov.Range("B6").Value = "Not to be localized"
$ov.Range("C6").Value = "Not to be localized"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = "In Translation"
$zh.Range("B3").Value = "In Translation"

$zh.Hyperlinks.Add($zh.Range("A4"), ($mdBase + $file1md), "", "", $file1md)
$zh.Range("B4").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("C4"), ($xlfZh + $file1xlfZh), "", "", $file1xlfZh)
$zh.Range("D4").Value = $handoffZh
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A5"), ($mdBase + $file2md), "", "", $file2md)
$zh.Range("B5").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("C5"), ($xlfZh + $file2xlfZh), "", "", $file2xlfZh)
$zh.Range("D5").Value = $handoffZh
$zh.Range("G5").Value = "0001-01-01 00:00:00"
$zh.Range("H5").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A6"), $cfgUrl, "", "", ".localization-config")
$zh.Range("B6").Value = "Not to be localized"
$zh.Range("D6").Value = "0001-01-01 00:00:00"
$zh.Range("G6").Value = "0001-01-01 00:00:00"
$zh.Range("H6").Value = "Ignored"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = "In Translation"
$de.Range("B3").Value = "In Translation"

$de.Hyperlinks.Add($de.Range("A4"), ($mdBase + $file1md), "", "", $file1md)
$de.Range("B4").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("C4"), ($xlfDe + $file1xlfDe), "", "", $file1xlfDe)
$de.Range("D4").Value = $handoffDe
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Include"

$de.Hyperlinks.Add($de.Range("A5"), ($mdBase + $file2md), "", "", $file2md)
$de.Range("B5").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("C5"), ($xlfDe + $file2xlfDe), "", "", $file2xlfDe)
$de.Range("D5").Value = $handoffDe
$de.Range("G5").Value = "0001-01-01 00:00:00"
$de.Range("H5").Value = "Include"

$de.Hyperlinks.Add($de.Range("A6"), $cfgUrl, "", "", ".localization-config")
$de.Range("B6").Value = "Not to be localized"
$de.Range("D6").Value = "0001-01-01 00:00:00"
$de.Range("G6").Value = "0001-01-01 00:00:00"
$de.Range("H6").Value = "Ignored"
